$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    'D2' = '46.749.97'
    'E2' = '  +4.12%  '
    'D3' = '2.256.12'
    'E3' = '  -0.17%  '
    'E4' = '  -0.10%  '
    'D5' = '300.41'
    'E5' = '  -0.68%  '
    'D6' = '100.30'
    'E6' = '  +6.15%  '
    'E7' = '  -0.53%  '
    'D8' = '0.999'
    'E8' = '  -0.04%  '
    'E9' = '  +0.77%  '
    'D10' = '35.43'
    'E10' = '  +4.16%  '
    'E11' = '  -1.07%  '
    'E12' = '  -0.28%  '
    'E13' = '  -0.75%  '
    'D14' = '2.602.26'
    'E14' = '  -0.19%  '
    'D15' = '2.257.40'
    'E15' = '  -0.28%  '
    'E16' = '  -0.27%  '
    'B17' = 'WrappedBTC'
    'C17' = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
    'D17' = '46.693.01'
    'E17' = '  +4.16%  '
    'B18' = 'Polygon'
    'C18' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
    'D18' = '0.793'
    'E18' = '  -0.46%  '
    'D19' = '12.79'
    'E19' = '  -0.29%  '
    'E20' = '  +0.43%  '
    'E21' = '  -2.85%  '
    'D22' = '65.34'
    'E22' = '  -0.22%  '
    'D23' = '249.45'
    'E23' = '  +4.96%  '
    'D24' = '2.82'
    'E24' = '  -2.45%  '
    'D25' = '1.00'
    'E25' = '  +0.26%  '
    'E26' = '  -0.82%  '
    'D27' = '42.49'
    'E27' = '  +2.64%  '
    'E28' = '  -0.66%  '
    'D29' = '9.66'
    'E29' = '  +1.00%  '
    'D30' = '19.89'
    'E30' = '  +1.87%  '
    'E31' = '  +8.75%  '
    'D32' = '146.06'
    'E32' = '  -4.11%  '
    'D33' = '5.40'
    'E33' = '  -2.59%  '
    'B34' = 'Hedera'
    'C34' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
    'D34' = '0.0767'
    'E34' = '  -2.57%  '
    'B35' = 'LidoDAOToken'
    'C35' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'D35' = '3.17'
    'E35' = '  +7.86%  '
    'D36' = '0.114'
    'E36' = '  +9.29%  '
    'E37' = '  -1.17%  '
    'D38' = '16.12'
    'E38' = '  +18.85%  '
    'E39' = '  -2.37%  '
    'D40' = '3.83'
    'E40' = '  -4.83%  '
    'D41' = '0.0297'
    'E41' = '  -4.76%  '
    'E42' = '  -1.64%  '
    'D43' = '0.998'
    'E43' = '  -0.21%  '
    'E44' = '  +1.78%  '
    'D45' = '1.819.86'
    'E45' = '  +4.25%  '
    'D46' = '90.70'
    'E46' = '  +19.87%  '
    'B47' = 'Algorand'
    'C47' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'D47' = '0.188'
    'E47' = '  -3.88%  '
    'B48' = 'ordi'
    'C48' = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
    'D48' = '71.99'
    'E48' = '  +2.52%  '
    'E49' = '  +2.48%  '
    'D50' = '93.73'
    'E50' = '  -2.29%  '
    'B51' = 'FraxShare'
    'C51' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
    'D51' = '7.75'
    'E51' = '  -0.57%  '
}

foreach ($ref in $updates.Keys) {
    $ws.Range($ref).Value = $updates[$ref]
}
